$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("week2")

# Set the "Actual time length to complete" for the Hand-in assignment row (row 18) to 5 hours (0.20833333333333334 as a day fraction)
$ws.Range("C18").Value = 0.20833333333333334

# Update the selection to C19 on the active sheet (week2), matching the post-edit state
$ws.Activate()
$ws.Range("C19").Select()

